$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A-C and E widened, D stays at the sheet default) ---
$ws.Columns.Item(1).ColumnWidth = 31.43   # A -> ~32.29
$ws.Columns.Item(2).ColumnWidth = 37.93   # B -> ~38.71
$ws.Columns.Item(3).ColumnWidth = 28.29   # C -> ~29.14
$ws.Columns.Item(5).ColumnWidth = 27.29   # E -> ~28.14

# --- Whole data body (existing + new rows) is (re)formatted as Text ---
$ws.Range("A2:D10").NumberFormat = "@"

# --- Row 5 ---
$ws.Cells.Item(5, 1).Value = "1898367659"
$ws.Cells.Item(5, 2).Value = "2024-12-01 10:23:25"
$ws.Cells.Item(5, 3).Value = "Нет"
$ws.Cells.Item(5, 4).Value = "Нет"

# --- Row 6 ---
$ws.Cells.Item(6, 1).Value = "1898367632"
$ws.Cells.Item(6, 2).Value = "2024-12-01 10:23:29"
$ws.Cells.Item(6, 3).Value = "Нет"
$ws.Cells.Item(6, 4).Value = "Нет"

# --- Row 7 ---
$ws.Cells.Item(7, 1).Value = "1898367594"
$ws.Cells.Item(7, 2).Value = "2024-12-10 15:00:00"
$ws.Cells.Item(7, 3).Value = "Да"
$ws.Cells.Item(7, 4).Value = "Да"

# --- Row 8 ---
$ws.Cells.Item(8, 1).Value = "1898367659"
$ws.Cells.Item(8, 2).Value = "2024-12-10 15:00:01"
$ws.Cells.Item(8, 3).Value = "Нет"
$ws.Cells.Item(8, 4).Value = "Нет"

# --- Row 9 ---
$ws.Cells.Item(9, 1).Value = "1898367660"
$ws.Cells.Item(9, 2).Value = "2024-12-15 23:02:02"
$ws.Cells.Item(9, 3).Value = "Да"
$ws.Cells.Item(9, 4).Value = "Да"

# --- Row 10 ---
$ws.Cells.Item(10, 1).Value = "1898367665"
$ws.Cells.Item(10, 2).Value = "2024-12-15 23:02:03"
$ws.Cells.Item(10, 3).Value = "Нет"
$ws.Cells.Item(10, 4).Value = "Нет"

# --- Selection ends up on B14, matching where the user left off ---
$ws.Range("B14").Select() | Out-Null
